$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.052176713943481
$ws.Range("B1").Value = 2.659528255462646
$ws.Range("C1").Value = 5.747958660125732
$ws.Range("D1").Value = 2.116188287734985
$ws.Range("E1").Value = 1.199236750602722
